$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" (column E) / "Valor Mora" (column F) rows 16-38 are
# reversed into reverse-chronological order (most recent period first),
# keeping every other column (B, C, D, G, H, I, J) untouched.

$periods = @(
    "2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312",
    "2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412"
)
$valores = @(
    247192,247192,247192,247192,247192,247192,247192,247192,247192,247192,247192,
    247192,247192,247192,247192,247192,247192,247192,247192,247192,247192,247192,189514
)

$n = $periods.Length
for ($i = 0; $i -lt $n; $i++) {
    $row = 16 + $i
    $srcIdx = $n - 1 - $i
    $ws.Range("E$row").Value = $periods[$srcIdx]
    $ws.Range("F$row").Value = $valores[$srcIdx]
}
